$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 get new (identical) random ID values
$ws.Range("D2").Value = "b2wRV2uP7N1iwZioFGrT"
$ws.Range("C2").Value = "XKZz3YI8pIlOahTDpCW2"
$ws.Range("E2").Value = "z3G8yCe31oIIfW9Hd2AD"
$ws.Range("F2").Value = "agzYZgUuhmN2UOEhBK8m"

$ws.Range("C3").Value = "XKZz3YI8pIlOahTDpCW2"
$ws.Range("D3").Value = "b2wRV2uP7N1iwZioFGrT"
$ws.Range("E3").Value = "z3G8yCe31oIIfW9Hd2AD"
$ws.Range("F3").Value = "agzYZgUuhmN2UOEhBK8m"

# Row 8 cleared out
$ws.Range("C8:F8").ClearContents()

# Row 9 gets new random ID values
$ws.Range("C9").Value = "K07DqGryxNUsknDgI999"
$ws.Range("D9").Value = "NTkDIAYnjMrRy5r6q2bK"
$ws.Range("E9").Value = "IEKuU0pORy6RuRPksnej"
$ws.Range("F9").Value = "HZiMBRZ2mvEApSy80Et1"

# Row 14 cleared out
$ws.Range("C14:F14").ClearContents()

# Row 15 gets new random ID values
$ws.Range("C15").Value = "KApyw3tfVKIqznJfqB6v"
$ws.Range("D15").Value = "pUriqnWv0JwDX3F03dwT"
$ws.Range("E15").Value = "X42LIVdTlVXxzETQeBT5"
$ws.Range("F15").Value = "oSDPwiYwofcwpRYpLNlr"

# Update the selected cell to match the author's final cursor position
$ws.Range("F15").Select()
